# Generate Report for Handback
# A new handback file (c48fc410-e840-42e8-bad0-39f5d71d25ad) was produced for
# row 7 of both the "zh-cn" and "de-de" localization-status sheets, but it is
# not the latest handback for that source file, so the status column records
# a warning while the handback file / datetime columns get filled in.

$wb = $excel.ActiveWorkbook

$currentHandbackUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/3b2799f731142694d680319c0498bbda8bab04ac/e2e/c48fc410-e840-42e8-bad0-39f5d71d25ad.md"
$latestHandbackUrl  = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/32180245582a6f46d3b6c8496420c44c2c553a1a/e2e/c48fc410-e840-42e8-bad0-39f5d71d25ad.md"
$notLatestMessage = "The version of handback file is not the latest, current: $currentHandbackUrl, latest: $latestHandbackUrl."
$handbackFileName = "c48fc410-e840-42e8-bad0-39f5d71d25ad.md"

# ---------------------------------------------------------------------------
# zh-cn sheet, row 7 (c48fc410-e840-42e8-bad0-39f5d71d25ad.md)
# ---------------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")

$wsZhCn.Range("I7").Value = $handbackFileName
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("I7"), $currentHandbackUrl, "", "", $handbackFileName)

$wsZhCn.Range("J7").Value = "c48fc410-e840-42e8-bad0-39f5d71d25ad.5bf15196ad42216c4b27acae77087e5e6453ab24.zh-cn.xlf"
$wsZhCn.Range("K7").Value = "2016-08-23 15:08:35"
$wsZhCn.Range("P7").Value = $notLatestMessage

# ---------------------------------------------------------------------------
# de-de sheet, row 7 (c48fc410-e840-42e8-bad0-39f5d71d25ad.md)
# ---------------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")

$wsDeDe.Range("I7").Value = $handbackFileName
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("I7"), $currentHandbackUrl, "", "", $handbackFileName)

$wsDeDe.Range("J7").Value = "c48fc410-e840-42e8-bad0-39f5d71d25ad.5bf15196ad42216c4b27acae77087e5e6453ab24.de-de.xlf"
$wsDeDe.Range("K7").Value = "2016-08-23 15:08:43"
$wsDeDe.Range("P7").Value = $notLatestMessage
